$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Range('D2').Value = '30.156.24'
$ws.Range('D3').Value = '1.853.98'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.13'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06398'
$ws.Range('D10').Value = '1.852.44'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07453'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.21'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.010'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '85.12'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6286'
$ws.Range('D16').Value = '30.122.38'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.37'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.71'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007318'
$ws.Range('D21').Value = '2.095.39'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.027'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.983'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.235'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.67'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1014'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.384'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.129'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.914'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04888'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.140'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7191'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.693'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01876'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.634'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9008'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.951'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '105.51'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9989'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4086'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.056'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '61.00'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1194'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.685'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.17'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.397'

# Update Volume(1h) (column E) values
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('E3').Value = '  -1.44%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  -3.09%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  -1.98%  '
$ws.Range('E8').Value = '  -3.51%  '
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('E10').Value = '  -1.55%  '
$ws.Range('E11').Value = '  -1.00%  '
$ws.Range('E12').Value = '  -2.45%  '
$ws.Range('E13').Value = '  -2.10%  '
$ws.Range('E14').Value = '  -4.05%  '
$ws.Range('E15').Value = '  -5.56%  '
$ws.Range('E16').Value = '  -1.43%  '
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('E19').Value = '  -4.94%  '
$ws.Range('E20').Value = '  -3.83%  '
$ws.Range('E21').Value = '  -3.71%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('E23').Value = '  -5.04%  '
$ws.Range('E24').Value = '  -3.23%  '
$ws.Range('E25').Value = '  -1.54%  '
$ws.Range('E26').Value = '  -1.24%  '
$ws.Range('E27').Value = '  -4.58%  '
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('E29').Value = '  +5.89%  '
$ws.Range('E30').Value = '  -2.71%  '
$ws.Range('E31').Value = '  -5.51%  '
$ws.Range('E32').Value = '  -3.18%  '
$ws.Range('E33').Value = '  -2.89%  '
$ws.Range('E34').Value = '  -5.62%  '
$ws.Range('E35').Value = '  -3.75%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('E38').Value = '  +1.31%  '
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('E40').Value = '  -1.80%  '
$ws.Range('E41').Value = '  -6.88%  '
$ws.Range('E42').Value = '  -0.93%  '
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('E44').Value = '  -4.64%  '
$ws.Range('E45').Value = '  -4.46%  '
$ws.Range('E46').Value = '  -4.87%  '
$ws.Range('E47').Value = '  -5.03%  '
$ws.Range('E48').Value = '  -7.15%  '
$ws.Range('E49').Value = '  -2.99%  '
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('E51').Value = '  -5.67%  '
